# Generate Report for Handoff
#
# The localization status for zh-cn / de-de moved from "In Translation" to
# "Ready for handoff", so:
#   - every "Status" cell (and the Overview rollup columns that mirror it)
#     gets the new text,
#   - the associated "generate"/"handoff" timestamps are refreshed, and
#   - the now-wider status columns are resized so the new text fits.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet: zh-cn / de-de status + the HO xliff generate date ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-09-02 00:48:42"

$overview.Columns.Item(5).ColumnWidth = 16.3333
$overview.Columns.Item(6).ColumnWidth = 16.3333

# --- zh-cn sheet: Status + Latest Handoff Datetime ------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-09-02 00:48:37"

$zhcn.Columns.Item(3).ColumnWidth = 16.3333

# --- de-de sheet: Status + Latest Handoff Datetime ------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-09-02 00:48:42"

$dede.Columns.Item(3).ColumnWidth = 16.3333
